$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A and append the new tickers right after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$newValues = @("LINK-USD", "MKR-USD", "GRT-USD")

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $lastRow + 1 + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
